# Auto-generated edit script applying the diff to 上海-漫展信息.xlsx
# Updates "F" (want-to-go count) and "G" (min price) columns across sheets,
# plus a full row-15 content refresh on sheet "全部类型" (event changed).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("展览")
$ws.Range("F5").Value = 2746
$ws.Range("F9").Value = 268
$ws.Range("F10").Value = 6224
$ws.Range("F14").Value = 4982
$ws.Range("F16").Value = 536
$ws.Range("F17").Value = 2593
$ws.Range("F18").Value = 1340
$ws.Range("F21").Value = 301
$ws.Range("F22").Value = 117
$ws.Range("F24").Value = 1047
$ws.Range("F25").Value = 229
$ws.Range("F26").Value = 387
$ws.Range("F28").Value = 1360
$ws.Range("F29").Value = 1020
$ws.Range("F30").Value = 2085
$ws.Range("F31").Value = 298
$ws.Range("F32").Value = 573
$ws.Range("F33").Value = 17
$ws.Range("F34").Value = 78
$ws.Range("F35").Value = 244
$ws.Range("F37").Value = 611
$ws.Range("F38").Value = 1038
$ws.Range("F41").Value = 14
$ws.Range("F42").Value = 282
$ws.Range("F43").Value = 2251
$ws.Range("F44").Value = 2536
$ws.Range("F49").Value = 85
$ws = $wb.Worksheets("演出")
$ws.Range("G4").Value = "不可售"
$ws.Range("F6").Value = 18
$ws.Range("F8").Value = 317
$ws.Range("F10").Value = 85
$ws.Range("F11").Value = 201
$ws.Range("F20").Value = 3
$ws.Range("F23").Value = 353
$ws = $wb.Worksheets("本地生活")
$ws.Range("F8").Value = 1462
$ws.Range("F10").Value = 2484
$ws.Range("F12").Value = 713
$ws = $wb.Worksheets("全部类型")
$ws.Range("F8").Value = 2746
$ws.Range("F10").Value = 1462
$ws.Range("F11").Value = 268
$ws.Range("F12").Value = 2484
$ws.Range("F13").Value = 6224
$ws.Range("B15").Value = "'2024-08-09"
$ws.Range("C15").Value = "上海·新井里美粉丝见面会"
$ws.Range("D15").Value = "西藏南路1号 上海大世界"
$ws.Range("E15").Value = "2024.08.09 16:30-08.11 16:30"
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 198
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=90332"
$ws.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202408/geUVjFXF1722842759283.jpeg"
$ws.Range("F17").Value = 4982
$ws.Range("F18").Value = 2593
$ws.Range("F21").Value = 301
$ws.Range("F22").Value = 117
$ws.Range("F24").Value = 1047
$ws.Range("F25").Value = 229
$ws.Range("F26").Value = 85
$ws.Range("F27").Value = 387
$ws.Range("F28").Value = 1360
$ws.Range("F29").Value = 1020
$ws.Range("F30").Value = 2085
$ws.Range("F31").Value = 298
$ws.Range("F32").Value = 573
$ws.Range("F33").Value = 244
$ws.Range("F36").Value = 1038
$ws.Range("F38").Value = 3
$ws.Range("F40").Value = 282
$ws.Range("F42").Value = 2251
$ws.Range("F43").Value = 2536
$ws.Range("F47").Value = 85
